$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.143.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.648.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3937"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.001"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.376"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.67%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.112"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001286"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.539"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.646.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06918"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.943"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.140.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.870"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.239"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "140.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.304"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.505"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.823.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08172"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.853"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02927"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9696"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.435"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7548"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6929"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.466"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.097"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08401"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.272"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.40%  "
